# Update literature review (to finish)
#
# Row 7 already carries the "in progress / highlighted" formatting
# (style indices: A=4, B=7, C/D/E=6) that rows 3 and 6 need to adopt,
# so copy that formatting down onto rows 3 and 6, then fill in the new
# Weight values for D3/D6. The dependent shared formulas in column E
# and the totals in row 35 recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the "done" row formatting (row 7) onto rows 3 and 6.
$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("A3:E3").PasteSpecial(-4122) | Out-Null

$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("A6:E6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Update the Weight (1-5) values for the two rows being finished.
$ws.Range("D3").Value = 100
$ws.Range("D6").Value = 100

# Move the active selection to D7, matching the saved view state.
$ws.Range("D7").Select() | Out-Null
